$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "IvánGarcía@beeckerco.com"
$ws.Range("A3").Value = "IvánGarcía@beeckerco.com"
$ws.Range("A4").Value = "AllissonFlores@beeckerco.com"
